$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 31 ("8/Server-Spiel" section) - fill in the previously blank data row.
# Copy formatting from the row above (row 30, same section block) so the
# cells pick up the existing "s=4 / s=3 / s=5" style indices used throughout
# the sheet instead of creating new ones.
# ---------------------------------------------------------------------------
$ws.Range("A30:H30").Copy()
$ws.Range("A31:H31").PasteSpecial(-4122)

$ws.Range("A31").Value = "8/Server-Spiel`nVersion: 1.1"
$ws.Range("A31").Characters(1, 14).Font.Bold = $true
$ws.Range("A31").Characters(15, 13).Font.Bold = $false

$ws.Range("B31").Value = "Funktional"
$ws.Range("C31").Value = "Länder die einem bestimmten `nSpieler nicht gehören, werden nicht`n""ausgegraut"", sondern in den Farben`ndes entsprechenden Besitzers ange-`nzeigt"
$ws.Range("D31").Value = "Spart Ressourcen,`nda jedem Client die`ngleiche Spielkarte`nangezeigt wird"
$ws.Range("E31").Value = "Bosin`n(10.10.2018)"
$ws.Range("F31").Value = "Kein"
$ws.Range("G31").Value = "6/Client-Init"
$ws.Range("H31").Value = "V1.1"
$ws.Range("H31").NumberFormat = "d-mmm"

$ws.Rows.Item(31).RowHeight = 90

# ---------------------------------------------------------------------------
# Row 19 ("5/Server-GUI" section) - fill in the previously blank data row.
# Copy formatting from the row above (row 18, same section block).
# ---------------------------------------------------------------------------
$ws.Range("A18:H18").Copy()
$ws.Range("A19:H19").PasteSpecial(-4122)

$ws.Range("A19").Value = "5/Server-GUI`nVersion: 1.1"
$ws.Range("A19").Characters(1, 12).Font.Bold = $true
$ws.Range("A19").Characters(13, 13).Font.Bold = $false

$ws.Range("B19").Value = "Funktional"
$ws.Range("C19").Value = "Erste Anzeige sind zwei Buttons zur`nAuswahl: neues Spiel und Spiel`nladen"
$ws.Range("D19").Value = "Die Unterscheidung`nist nötig, da nur `nbei einem neuen `nSpiel die Spieler-`nanzahl festgelegt`nwerden muss."
$ws.Range("E19").Value = "Kunde`n(10.10.2018)"
$ws.Range("F19").Value = "Kein"
$ws.Range("G19").Value = "4/Server-Init"
$ws.Range("H19").Value = "V1.1"

$ws.Rows.Item(19).RowHeight = 105

# ---------------------------------------------------------------------------
# Update the active selection to match the edited area.
# ---------------------------------------------------------------------------
[void]$ws.Range("D26").Select()
